$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 6, shifting existing rows 6-17 down to 7-18.
$ws.Rows("6:6").Insert()

# New row 6 uses the same "wrap text / vertical-center" formatting as the
# rest of the data rows (cellXf s="2").
$ws.Range("A6:I6").WrapText = $true
$ws.Range("A6:I6").VerticalAlignment = -4108
$ws.Rows("6:6").RowHeight = 105

# Fill in the new row 6 content.
$ws.Range("A6").Value = "Aus aus"
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = "Alaska"
$ws.Range("I6").Value = "Smith, 1920. Bears on the coast. Jr. Chilly Waters. 0:0 pp0-40."

# Update the view: scroll back to the top and select H6 (matches target).
$ws.Range("H6").Select()
